# Apply "added background image customization" edit.
$wb = $excel.ActiveWorkbook

$debugWs = $wb.Worksheets.Item("Debug Questions")
$quizWs  = $wb.Worksheets.Item("Example")

# Rename the "Example" sheet to "Science Quiz" - the feature moves from the
# debug/example sheet onto the real quiz sheet.
$quizWs.Name = "Science Quiz"

# On the Debug Questions sheet, rename the header labels for the Hint/Background
# image columns (H1/I1) from the underscore form to a space form, and remove the
# sample/example image-url row (H2/I2) that was only there for illustration.
$debugWs.Range("H1").Value = "Hint Image"
$debugWs.Range("I1").Value = "Background Image"
$debugWs.Range("H2:I2").ClearContents()

# Add the same Hint Image / Background Image header columns to the quiz sheet,
# actually enabling the background-image customization feature there.
$quizWs.Range("H1").Value = "Hint Image"
$quizWs.Range("I1").Value = "Background Image"
$quizWs.Range("H1:I1").Style = "Normal"

# Update the active sheet / selection so the workbook opens on the quiz sheet.
$debugWs.Range("D28").Select()
$quizWs.Range("C17").Select()
$quizWs.Activate()
